# Update crypto price/volume data per the GitHub Actions scraper run (Wed May  3 14:21:36 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column whose new values look like plain numbers need an explicit
# text format first, otherwise Excel auto-converts them (e.g. "21.50" -> 21.5, losing the
# trailing zero / grouping dots) instead of keeping the scraped string verbatim.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values (Coin/Link/Price/Volume(1h)) row by row.
$ws.Range('D2').Value = '28.473.68'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '1.864.35'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '324.36'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.4548'
$ws.Range('E7').Value = '  -2.09%  '
$ws.Range('D8').Value = '0.3826'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').Value = '0.07815'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').Value = '0.9864'
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').Value = '21.50'
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('D12').Value = '1.859.60'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = '6.893'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').Value = '5.611'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').Value = '0.06908'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '86.64'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '0.000009924'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = '16.61'
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = '28.486.55'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').Value = '5.236'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = '10.87'
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('D24').Value = '2.097'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '2.081.82'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').Value = '153.68'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').Value = '19.07'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = '5.654'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '117.26'
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = '1.916'
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('D31').Value = '0.09254'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').Value = '0.9030'
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('D33').Value = '5.243'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('D35').Value = '3.294'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').Value = '0.05685'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('D39').Value = '7.643'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('D40').Value = '0.5541'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('D41').Value = '0.1766'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').Value = '9.608'
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('D43').Value = '0.07082'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').Value = '11.50'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('D45').Value = '0.5226'
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('D46').Value = '1.131'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('D47').Value = '2.102'
$ws.Range('E47').Value = '  -1.95%  '
$ws.Range('D48').Value = '1.806'
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('D49').Value = '111.76'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('D50').Value = '2.427'
$ws.Range('E50').Value = '  +4.32%  '
$ws.Range('D51').Value = '1.006'
$ws.Range('E51').Value = '  -0.01%  '
